# Restore cell C10 on the active ("Rules") sheet to its prior saved value.
# Per the commit diff, cell C10 (the "From" value for rule R30) changes
# from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
